# Fruta / hortaliza, semanal
# A new weekly price record for "Ajo" (Femacal de La Calera) is inserted
# as a new row 211 in the data table; every subsequent record shifts down
# by one row (old row 211 -> new row 212, ..., old row 237 -> new row 238).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 211..237 down to 212..238, leaving a blank row 211 behind.
$ws.Rows("211:211").Insert()

# Populate the newly inserted row 211 with the new record's data.
$ws.Cells.Item(211, 1).Value2 = 3
$ws.Cells.Item(211, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(211, 3).Value2 = "Coquimbo"
$ws.Cells.Item(211, 4).Value2 = 44491
$ws.Cells.Item(211, 5).Value2 = 5
$ws.Cells.Item(211, 6).Value2 = 100112003
$ws.Cells.Item(211, 7).Value2 = "Ajo"
$ws.Cells.Item(211, 8).Value2 = "Chino"
$ws.Cells.Item(211, 9).Value2 = "Primera"
$ws.Cells.Item(211, 10).Value2 = 73
$ws.Cells.Item(211, 11).Value2 = 15500
$ws.Cells.Item(211, 12).Value2 = 16000
$ws.Cells.Item(211, 13).Value2 = 15760
$ws.Cells.Item(211, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(211, 15).Value2 = "China"
$ws.Cells.Item(211, 16).Value2 = 1576
$ws.Cells.Item(211, 17).Value2 = 10
$ws.Cells.Item(211, 18).Value2 = "Hortaliza"
